$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new "Price" values look like plain numbers (single decimal
# point, e.g. "585.27"), which Excel would otherwise auto-convert to a numeric
# cell. The source data stores these as literal text (inlineStr), so force the
# Text format before assigning, then restore the Normal cell style so no stray
# number-format style is left behind on the cell.
$textCells = @('D5','D6','D8','D9','D10','D11','D13','D20','D21','D23','D25','D28','D30','D34','D35','D37','D38','D39','D41','D42','D45','D46','D48','D50')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = '@'
}

$ws.Range('D2').Value = '63.306.78'
$ws.Range('E2').Value = '  +2.82%  '
$ws.Range('D3').Value = '3.488.28'
$ws.Range('E3').Value = '  +2.86%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '585.27'
$ws.Range('E5').Value = '  +1.42%  '
$ws.Range('D6').Value = '148.14'
$ws.Range('E6').Value = '  +5.51%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.481'
$ws.Range('E8').Value = '  +1.22%  '
$ws.Range('D9').Value = '7.71'
$ws.Range('E9').Value = '  +0.73%  '
$ws.Range('D10').Value = '0.127'
$ws.Range('E10').Value = '  +3.19%  '
$ws.Range('D11').Value = '0.399'
$ws.Range('E11').Value = '  +2.71%  '
$ws.Range('D12').Value = '4.084.87'
$ws.Range('E12').Value = '  +2.89%  '
$ws.Range('D13').Value = '29.83'
$ws.Range('E13').Value = '  +4.78%  '
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').Value = '3.488.87'
$ws.Range('E15').Value = '  +2.64%  '
$ws.Range('E16').Value = '  +2.18%  '
$ws.Range('D17').Value = '63.316.30'
$ws.Range('E17').Value = '  +2.87%  '
$ws.Range('E18').Value = '  +2.25%  '
$ws.Range('E19').Value = '  +5.35%  '
$ws.Range('D20').Value = '9.37'
$ws.Range('E20').Value = '  +4.37%  '
$ws.Range('D21').Value = '390.91'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('E22').Value = '  +1.63%  '
$ws.Range('D23').Value = '75.23'
$ws.Range('E23').Value = '  +0.26%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').Value = '0.0000118'
$ws.Range('E25').Value = '  +5.33%  '
$ws.Range('D26').Value = '3.629.06'
$ws.Range('E26').Value = '  +2.96%  '
$ws.Range('E27').Value = '  -4.37%  '
$ws.Range('D28').Value = '7.84'
$ws.Range('E28').Value = '  +7.99%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = '8.32'
$ws.Range('E30').Value = '  +3.43%  '
$ws.Range('E31').Value = '  +6.99%  '
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').Value = '23.88'
$ws.Range('E34').Value = '  +1.73%  '
$ws.Range('D35').Value = '5.37'
$ws.Range('E35').Value = '  +6.46%  '
$ws.Range('E36').Value = '  +2.93%  '
$ws.Range('D37').Value = '32.22'
$ws.Range('E37').Value = '  +24.26%  '
$ws.Range('D38').Value = '171.41'
$ws.Range('E38').Value = '  +2.10%  '
$ws.Range('D39').Value = '1.58'
$ws.Range('E39').Value = '  +7.21%  '
$ws.Range('D40').Value = '3.524.32'
$ws.Range('E40').Value = '  +2.87%  '
$ws.Range('D41').Value = '0.0771'
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').Value = '0.810'
$ws.Range('E42').Value = '  +4.20%  '
$ws.Range('E43').Value = '  +1.66%  '
$ws.Range('E44').Value = '  +4.56%  '
$ws.Range('D45').Value = '42.45'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').Value = '1.22'
$ws.Range('E46').Value = '  +7.19%  '
$ws.Range('D47').Value = '2.619.47'
$ws.Range('E47').Value = '  +6.73%  '
$ws.Range('D48').Value = '23.73'
$ws.Range('E48').Value = '  +5.17%  '
$ws.Range('E49').Value = '  +13.67%  '
$ws.Range('D50').Value = '6.79'
$ws.Range('E50').Value = '  +1.35%  '
$ws.Range('E51').Value = '  +3.54%  '

foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = 'Normal'
}
